# Outstandings.xlsx — "Add files via upload" edit
#
# Net effect (per the target diff): the purchase-side entry for
# "Sanyo and Sanyo" / Chq no 089933 (Sr. No 5, row 10, with its blank
# spacer row 11) is removed from the "Purchase 22-23" sheet. Every
# following entry shifts up two rows and its "Sr. No" is renumbered
# down by one. The now-unused "Chq no 089933" shared string drops out
# on save. The two sheets' remembered selections also move (an
# artifact of where the user's cursor ended up after the edit).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Purchase 22-23"
$ws2 = $wb.Worksheets.Item(2)   # "Sale 22-23"

# Remove the data row (10) and its trailing blank spacer row (11) for
# the deleted "Sanyo and Sanyo" / Chq no 089933 entry. Everything below
# shifts up by two rows, formulas (e.g. =E16+E17) retarget automatically.
$ws1.Rows("10:11").Delete()

# Renumber the "Sr. No" column for the remaining entries (each drops by 1).
$ws1.Range("A10").Value = 5
$ws1.Range("A12").Value = 6
$ws1.Range("A14").Value = 7
$ws1.Range("A17").Value = 8
$ws1.Range("A19").Value = 9

# Restore the recorded cursor positions on both sheets. Select sheet2's
# cell first so the final ".Select()" on sheet1 leaves "Purchase 22-23"
# as the active/selected tab, matching the original tabSelected state.
$ws2.Range("D29").Select()
$ws1.Range("D27").Select()
